$d = $word.ActiveDocument

$d.Content.Find.Execute("Play Free Mystery Joker Slot Game Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Mystery Joker Free: Classic Fruit Slot Game", 2)

$d.Content.Find.Execute("Exciting bonus features and big win opportunities", $true, $false, $false, $false, $false, $true, 1, $false, "Simple gameplay system", 2)

$d.Content.Find.Execute("High volatility for substantial wins", $true, $false, $false, $false, $false, $true, 1, $false, "Exciting bonus features", 2)

$d.Content.Find.Execute("Well-designed graphics from Play'n Go", $true, $false, $false, $false, $false, $true, 1, $false, "Well-designed graphics", 2)

$d.Content.Find.Execute("Autoplay options for convenience", $true, $false, $false, $false, $false, $true, 1, $false, "Reputable software provider", 2)

$d.Content.Find.Execute("Low RTP", $true, $false, $false, $false, $false, $true, 1, $false, "High volatility", 2)

$d.Content.Find.Execute("Limited number of paylines", $true, $false, $false, $false, $false, $true, 1, $false, "RTP not guaranteed", 2)

$d.Content.Find.Execute("Read our review of Mystery Joker, a slot game with classic fruit themes and exciting features. Play for free and win big with Autoplay and free spins.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Mystery Joker, a classic fruit slot game with exciting features. Play for free now.", 2)
